$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows above the current row 30 (VERIFY_TEXT_PRESENT row),
# shifting it and subsequent rows down to make room for two new WAIT rows.
$ws.Rows.Item(30).Resize(2).Insert()

# New row 30: WAIT
$ws.Cells.Item(30, 1).Style = $ws.Cells.Item(32, 1).Style
$ws.Cells.Item(30, 2).Style = $ws.Cells.Item(32, 2).Style
$ws.Cells.Item(30, 3).Style = $ws.Cells.Item(32, 3).Style
$ws.Cells.Item(30, 4).Style = $ws.Cells.Item(32, 4).Style
$ws.Cells.Item(30, 5).Style = $ws.Cells.Item(32, 5).Style
$ws.Cells.Item(30, 2).Value = "WAIT"

# New row 31: WAIT
$ws.Cells.Item(31, 1).Style = $ws.Cells.Item(32, 1).Style
$ws.Cells.Item(31, 2).Style = $ws.Cells.Item(32, 2).Style
$ws.Cells.Item(31, 3).Style = $ws.Cells.Item(32, 3).Style
$ws.Cells.Item(31, 4).Style = $ws.Cells.Item(32, 4).Style
$ws.Cells.Item(31, 5).Style = $ws.Cells.Item(32, 5).Style
$ws.Cells.Item(31, 2).Value = "WAIT"

# Update selection to match the recorded edit point
$ws.Range("B31").Select()
